$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D20").Value = "[파이썬 간단한 게임 만들기] 9. 오목 아니고 4목"
$ws.Range("E20").Value = "https://ai-creator.tistory.com/537"

$ws.Range("D46").Value = "[SK바이오팜] 2021년 04월, 생물정보학(Bioinformatics 채용), Digital Healthcare 경력직 채용(신약개발 Big-data AI 모델 개발 및 분석)"
$ws.Range("E46").Value = "https://bioinformaticsandme.tistory.com/389"

$ws.Range("D51").Value = "첫번째, 두번째 vs 첫 번째, 두 번째"
$ws.Range("E51").Value = "https://bskyvision.com/1162"
